# Swap "email, System" -> "System, email" in the "Recorded By" column (G)
# for rows recorded by dnasr281@gmail.com or admin@admin.com (but not the
# backup@backdoor.com rows, which keep their original order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string]) {
        $parts = $val -split ', '

        if ($parts.Count -eq 2 -and $parts[1] -eq 'System' -and $parts[0] -ne 'backup@backdoor.com') {
            $cell.Value = 'System, ' + $parts[0]
        }
    }
}
